$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 09:16"

$data = @(
  ,@(4, "China", 81008,41,71740,6013,1927,7,3255)
  ,@(5, "Italia", 47021,0,5129,37860,2655,0,4032)
  ,@(6, "España", 21571,0,1588,18890,939,0,1093)
  ,@(7, "Alemania", 19848,0,180,19600,2,0,68)
  ,@(8, "Estados Unidos", 19774,391,147,19352,64,19,275)
  ,@(9, "Iran", 19644,0,6745,11466,0,0,1433)
  ,@(10, "Francia", 12612,0,1587,10575,1297,0,450)
  ,@(11, "Corea del Sur", 8799,147,2612,6085,59,8,102)
  ,@(12, "Suiza", 5615,0,15,5544,141,0,56)
  ,@(13, "Reino Unido", 3983,0,65,3741,20,0,177)
  ,@(14, "Paises Bajos", 2994,0,2,2886,210,0,106)
  ,@(15, "Austria", 2666,17,9,2651,14,0,6)
  ,@(16, "Belgica", 2257,0,204,2016,164,0,37)
  ,@(17, "Noruega", 1975,16,1,1967,27,0,7)
  ,@(18, "Suecia", 1639,0,16,1607,21,0,16)
  ,@(19, "Dinamarca", 1255,0,1,1245,37,0,9)
  ,@(20, "Malasia", 1183,153,114,1065,26,1,4)
  ,@(21, "Canada", 1087,0,14,1061,1,0,12)
  ,@(22, "Australia", 1068,140,46,1015,2,0,7)
  ,@(23, "Portugal", 1020,0,5,1009,26,0,6)
  ,@(24, "Japon", 1007,0,215,757,50,0,35)
  ,@(25, "Brasil", 977,7,2,964,18,0,11)
  ,@(26, "Chequia", 904,71,4,900,6,0,0)
  ,@(27, "Crucero", 712,0,567,137,15,0,8)
  ,@(28, "Israel", 705,0,15,689,10,0,1)
  ,@(29, "Irlanda", 683,0,5,675,6,0,3)
  ,@(30, "Turquia", 670,0,0,661,0,0,9)
  ,@(31, "Pakistan", 588,87,13,571,0,1,4)
  ,@(32, "Grecia", 495,0,19,466,20,0,10)
  ,@(33, "Luxemburgo", 484,0,6,473,1,0,5)
  ,@(34, "Catar", 470,0,10,460,6,0,0)
  ,@(35, "Finlandia", 450,0,10,440,2,0,0)
  ,@(36, "Polonia", 439,14,13,421,3,0,5)
  ,@(37, "Chile", 434,0,6,428,7,0,0)
  ,@(38, "Ecuador", 426,0,3,416,2,0,7)
  ,@(39, "Tailandia", 411,89,42,368,1,0,1)
  ,@(40, "Islandia", 409,0,5,404,1,0,0)
  ,@(41, "Singapur", 385,0,131,252,14,2,2)
  ,@(42, "Indonesia", 369,0,17,320,0,0,32)
  ,@(43, "Arabia Saudita", 344,0,8,336,0,0,0)
  ,@(44, "Eslovenia", 341,0,0,340,9,0,1)
  ,@(45, "Rumania", 308,0,31,277,11,0,0)
  ,@(46, "Barein", 298,0,125,172,4,0,1)
  ,@(47, "Egipto", 285,0,42,235,0,0,8)
  ,@(48, "Estonia", 283,0,1,282,1,0,0)
  ,@(49, "India", 275,26,23,247,0,0,5)
  ,@(50, "Peru", 263,0,1,258,5,0,4)
  ,@(51, "Filipinas", 262,32,8,236,1,0,18)
  ,@(52, "Hong Kong", 256,0,98,154,4,0,4)
  ,@(53, "Rusia", 253,0,12,240,0,0,1)
  ,@(54, "Irak", 208,0,49,142,0,0,17)
  ,@(55, "Mexico", 203,39,4,197,1,1,2)
  ,@(56, "Sudafrica", 202,0,0,202,0,0,0)
  ,@(57, "Panama", 200,0,1,198,7,0,1)
  ,@(58, "Libano", 177,0,4,169,3,0,4)
  ,@(59, "Armenia", 160,24,1,159,2,0,0)
  ,@(60, "Kuwait", 159,0,22,137,5,0,0)
  ,@(61, "Colombia", 158,13,1,157,0,0,0)
  ,@(62, "Argentina", 158,0,3,152,0,0,3)
  ,@(63, "Croacia", 157,27,5,151,0,0,1)
  ,@(64, "Taiwan", 153,18,28,123,0,0,2)
  ,@(65, "San Marino", 151,0,4,133,12,0,14)
  ,@(66, "Serbia", 149,14,2,146,4,0,1)
  ,@(67, "Bulgaria", 142,15,3,136,3,0,3)
  ,@(68, "Emiratos Arabes Unidos", 140,0,31,107,2,0,2)
  ,@(69, "Eslovaquia", 137,0,0,137,2,0,0)
  ,@(70, "Letonia", 124,13,1,123,0,0,0)
  ,@(71, "Costa Rica", 113,0,2,109,2,0,2)
  ,@(72, "Uruguay", 110,0,0,110,0,0,0)
  ,@(73, "Hungria", 103,18,7,92,6,0,4)
  ,@(74, "Argelia", 94,0,32,51,0,0,11)
  ,@(75, "Vietnam", 91,0,17,74,0,0,0)
  ,@(76, "Bosnia y Herzegovina", 90,1,2,88,1,0,0)
  ,@(77, "Marruecos", 86,0,2,81,1,0,3)
  ,@(78, "Jordania", 85,0,1,84,0,0,0)
  ,@(79, "Islas Feroe", 80,0,3,77,0,0,0)
  ,@(80, "Brunei", 78,0,1,77,2,0,0)
  ,@(81, "Republica de Macedonia", 76,0,1,75,1,0,0)
  ,@(82, "Republica de Chipre", 75,0,0,75,1,0,0)
  ,@(83, "Principado de Andorra", 75,0,1,74,2,0,0)
  ,@(84, "Sri Lanka", 73,0,3,70,0,0,0)
  ,@(85, "Republica Dominicana", 72,0,0,70,0,0,2)
  ,@(86, "Albania", 70,0,2,66,2,0,2)
  ,@(87, "Lituania", 69,0,1,67,1,0,1)
  ,@(88, "Bielorrusia", 69,0,15,54,0,0,0)
  ,@(89, "Moldavia", 66,0,1,64,3,0,1)
  ,@(90, "Venezuela", 65,0,1,64,0,0,0)
  ,@(91, "Malta", 64,0,2,62,1,0,0)
  ,@(92, "Tunez", 54,0,1,52,7,0,1)
  ,@(93, "Kazajistan", 53,1,0,53,0,0,0)
  ,@(94, "Nueva Zelanda", 52,0,0,52,0,0,0)
  ,@(95, "Oman", 52,4,13,39,0,0,0)
  ,@(96, "Estado de Palestina", 52,4,17,35,0,0,0)
  ,@(97, "Camboya", 51,0,1,50,0,0,0)
  ,@(98, "Guadalupe", 51,0,0,50,4,0,1)
  ,@(99, "Georgia", 47,3,1,46,1,0,0)
  ,@(100, "Senegal", 47,0,5,42,0,0,0)
  ,@(101, "Azerbaiyan", 44,0,7,36,0,0,1)
  ,@(102, "Ucrania", 41,0,1,37,0,0,3)
  ,@(103, "Burkina Faso", 40,0,4,35,0,0,1)
  ,@(104, "Reunion", 38,0,0,38,0,0,0)
  ,@(105, "Liechtenstein", 37,9,0,37,0,0,0)
  ,@(106, "Uzbekistan", 37,4,0,37,0,0,0)
  ,@(107, "Martinica", 32,0,0,31,7,0,1)
  ,@(108, "Camerun", 27,0,2,25,0,0,0)
  ,@(109, "Honduras", 24,0,0,24,0,0,0)
  ,@(110, "Afganistan", 24,0,1,23,0,0,0)
  ,@(111, "Consejo Danes para los Refugiados", 23,5,0,22,0,1,1)
  ,@(112, "Cuba", 21,0,0,20,0,0,1)
  ,@(113, "Banglades", 20,0,3,16,1,0,1)
  ,@(114, "Bolivia", 19,3,0,19,0,0,0)
  ,@(115, "Jamaica", 19,0,2,16,0,0,1)
  ,@(116, "Paraguay", 18,0,0,18,1,0,0)
  ,@(117, "Ruanda", 17,0,0,17,0,0,0)
  ,@(118, "Macao", 17,0,10,7,0,0,0)
  ,@(119, "Ghana", 16,0,0,16,0,0,0)
  ,@(120, "Guayana Francesa", 15,0,0,15,0,0,0)
  ,@(121, "Polinesia Francesa", 15,4,0,15,0,0,0)
  ,@(122, "Guyana", 15,0,0,14,0,0,1)
  ,@(123, "Puerto Rico", 14,0,0,14,0,0,0)
  ,@(124, "Guam", 14,0,0,14,0,0,0)
  ,@(125, "Montenegro", 14,0,0,14,0,0,0)
  ,@(126, "Mauricio", 14,2,0,13,0,1,1)
  ,@(127, "Costa de Marfil", 14,5,1,13,0,0,0)
  ,@(128, "Maldivas", 13,0,2,11,0,0,0)
  ,@(129, "Kirguistan", 12,6,0,12,0,0,0)
  ,@(130, "Guatemala", 12,0,0,11,0,0,1)
  ,@(131, "Nigeria", 12,0,1,11,0,0,0)
  ,@(132, "Monaco", 11,0,0,11,0,0,0)
  ,@(133, "Mongolia", 10,4,0,10,0,0,0)
  ,@(134, "Gibraltar", 10,0,2,8,0,0,0)
  ,@(135, "Etiopia", 9,0,0,9,0,0,0)
  ,@(136, "Trinidad yTobago", 9,0,0,9,0,0,0)
  ,@(137, "Togo", 9,0,0,9,0,0,0)
  ,@(138, "Kenia", 7,0,0,7,0,0,0)
  ,@(139, "Mayotte", 7,0,0,7,0,0,0)
  ,@(140, "Seychelles", 7,0,0,7,0,0,0)
  ,@(141, "Barbados", 6,0,0,6,0,0,0)
  ,@(142, "Guinea Ecuatorial", 6,0,0,6,0,0,0)
  ,@(143, "Islas Virgenes de los Estados Unidos", 6,3,0,6,0,0,0)
  ,@(144, "Tanzania", 6,0,0,6,0,0,0)
  ,@(145, "Aruba", 5,0,1,4,0,0,0)
  ,@(146, "Surinam", 4,0,0,4,0,0,0)
  ,@(147, "Bahamas", 4,0,0,4,0,0,0)
  ,@(148, "San Martin (Parte Francesa)", 4,0,0,4,0,0,0)
  ,@(149, "Gabon", 4,0,0,3,0,0,1)
  ,@(150, "Namibia", 3,0,0,3,0,0,0)
  ,@(151, "Republica de Africa Central", 3,0,0,3,0,0,0)
  ,@(152, "San Bartolome", 3,0,0,3,0,0,0)
  ,@(153, "El Salvador", 3,2,0,3,0,0,0)
  ,@(154, "Congo", 3,0,0,3,0,0,0)
  ,@(155, "Madagascar", 3,0,0,3,0,0,0)
  ,@(156, "Islas Caimanes", 3,0,0,2,0,0,1)
  ,@(157, "Curazao", 3,0,0,2,0,0,1)
  ,@(158, "Zambia", 2,0,0,2,0,0,0)
  ,@(159, "Benin", 2,0,0,2,0,0,0)
  ,@(160, "Bermudas", 2,0,0,2,0,0,0)
  ,@(161, "Isla de Man", 2,0,0,2,0,0,0)
  ,@(162, "Guinea", 2,0,0,2,0,0,0)
  ,@(163, "Groenlandia", 2,0,0,2,0,0,0)
  ,@(164, "Fiyi", 2,1,0,2,0,0,0)
  ,@(165, "Santa Lucia", 2,0,0,2,0,0,0)
  ,@(166, "Nueva Caledonia", 2,0,0,2,0,0,0)
  ,@(167, "Nicaragua", 2,0,0,2,0,0,0)
  ,@(168, "Mauritania", 2,0,0,2,0,0,0)
  ,@(169, "Butan", 2,0,0,2,0,0,0)
  ,@(170, "Liberia", 2,0,0,2,0,0,0)
  ,@(171, "Haiti", 2,0,0,2,0,0,0)
  ,@(172, "Sudan", 2,0,0,1,0,0,1)
  ,@(173, "Montserrat", 1,0,0,1,0,0,0)
  ,@(174, "Republica del Chad", 1,0,0,1,0,0,0)
  ,@(175, "Timor Oriental", 1,1,0,1,0,0,0)
  ,@(176, "Niger", 1,0,0,1,0,0,0)
  ,@(177, "San Martin (Parte Holandesa)", 1,0,0,1,0,0,0)
  ,@(178, "Zimbabue", 1,0,0,1,0,0,0)
  ,@(179, "Republica de Yibuti", 1,0,0,1,0,0,0)
  ,@(180, "Suazilandia", 1,0,0,1,0,0,0)
  ,@(181, "Gambia", 1,0,0,1,0,0,0)
  ,@(182, "Cabo Verde", 1,0,0,1,0,0,0)
  ,@(183, "Santa Sede", 1,0,0,1,0,0,0)
  ,@(184, "Angola", 1,0,0,1,0,0,0)
  ,@(185, "Somalia", 1,0,0,1,0,0,0)
  ,@(186, "Papua Nueva Guinea", 1,0,0,1,0,0,0)
  ,@(187, "San Vicente y las Granadinas", 1,0,0,1,0,0,0)
  ,@(188, "Antigua y Barbuda", 1,0,0,1,0,0,0)
  ,@(189, "Nepal", 1,0,1,0,0,0,0)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
  $ws.Cells.Item($r, 6).Value = $row[6]
  $ws.Cells.Item($r, 7).Value = $row[7]
  $ws.Cells.Item($r, 8).Value = $row[8]
}
